# Update Name of Algo
# Apply the numeric value changes described by the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = 8.458200000000005
$ws.Range("B12").Value = 6.152400000000001
$ws.Range("E13").Value = 12.2307
$ws.Range("B18").Value = 4.737800000000004

$wb.Save()
